$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 159
$ws.Range("I2").Value = 403
$ws.Range("J2").Value = 1798
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 458
$ws.Range("M2").Value = 28
$ws.Range("N2").Value = 337
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 9
$ws.Range("R2").Value = 31
$ws.Range("S2").Value = 169
$ws.Range("T2").Value = 266
$ws.Range("V2").Value = 2664
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2693
$ws.Range("Z2").Value = 43
$ws.Range("AA2").Value = 13
